# Replace the literal "<br/>" line-break markers inside the Sector/Size
# threshold descriptions (rows 20-23, columns B-D on the "Summary" sheet)
# with real newline characters, so the text wraps as actual multi-line
# cell content instead of containing a literal HTML tag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "<20 Heavy Industry. `n<5 Wholesale `n<10 Retail, Accmmdtn., RstrnTurnover, Infrmtn., Sftwr., Tenancy.`n<20 Transport, Warehouse, Postal`n<100 Property management `n<10 Other Industries "
$ws.Range("C20").Value = "<3Millionlion CnstrcTurnover`n<20Millionlion Real Estate`n<1Millionlion Leasing and Commercial Services"
$ws.Range("D20").Value = "<3Millionlion Heavy Industry. CnstrcTurnover`n<10Millionlion Wholesale `n<1Millionlion Retail, Warehouse, Postal, Accmmdtn., RstrnTurnover, Infrmtn., Real State, Tenancy.`n<0.5Millionlion Sftwr., `n<2Millionlion Transport,`n<5Millionlion Property management"

$ws.Range("B21").Value = ">=20 Heavy Industry. `n>=5 Wholesale `n>=10 Retail, Accmmdtn., RstrnTurnover, Infrmtn., Sftwr., Tenancy.`n>=20 Transport, Warehouse, Postal`n>=100 Property management `n>=10 Other Industries "
$ws.Range("C21").Value = "<50Millionlion CnstrcTurnover`n<50Millionlion Real Estate`n< 80Millionlion Leasing and Commercial Services"
$ws.Range("D21").Value = ">=3Millionlion Heavy Industry. CnstrcTurnover`n>=10Millionlion Wholesale `n>=1Millionlion Retail, Warehouse, Postal, Accmmdtn., RstrnTurnover, Infrmtn., Real State, Tenancy.`n>=0.5Millionlion Sftwr., `n>=2Millionlion Transport,`n>=5Millionlion Property management"

$ws.Range("B22").Value = ">=300 Heavy Industry, Transport, Postal, Property management `n>=20 Wholesale `n>=50 Retail, `n>=100 Warehouse, Accmmdtn., RstrnTurnover, Infrmtn., Sftwr., Tenancy.`n>=100 Other Industries "
$ws.Range("C22").Value = "<800Millionlion CnstrcTurnover`n<100Millionlion Real Estate`n<120Millionlion Leasing and Commercial Services"
$ws.Range("D22").Value = ">=20Millionlion Heavy Industry, Postal, Accmmdtn., RstrnTurnover, `n>=60Millionlion CnstrcTurnover`n>=50Millionlion Wholesale `n>=5Millionlion Retail, `n>=10Millionlion Warehouse, Infrmtn., Sftwr., Real State, Property management`n>=80Millionlion Tenancy.`n>=30Millionlion Transport,"

$ws.Range("B23").Value = ">=1000 Heavy Industry, Transport, Postal, Property management `n>=200 Wholesale, Warehouse `n>=300 Retail, Accmmdtn., RstrnTurnover, Sftwr., Tenancy.`n>=2000 Infrmtn.,`n>=300 Other Industries "
$ws.Range("C23").Value = ">800Millionlion CnstrcTurnover`n>100Millionlion Real Estate`n>120Millionlion Leasing and Commercial Services"
$ws.Range("D23").Value = ">=400Millionlion Heavy Industry, Wholesale, `n>=100Millionlion Accmmdtn., RstrnTurnover, Sftwr.`n>=800Millionlion CnstrcTurnover`n>=200Millionlion Retail, `n>=1Billionlion  Infrmtn., `n>=2Billionlion Real State`n>=50Millionlion Property management`n>=1.2Billionlion Tenancy.`n>=300Millionlion Transport, Warehouse,"
